$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(39400, 0.7960900427496256),
    @(39583, 2.427463180990628),
    @(39765, -0.2156301208428602),
    @(39948, -9.679894544238294),
    @(40130, -0.8195522748085864),
    @(40310, 2.274906658561562),
    @(40494, 7.912074126222677),
    @(40676, 1.840273951972307),
    @(40862, 2.530590017523977),
    @(41044, 1.726328558848905),
    @(41228, 1.446992732417414),
    @(41409, -1.803419669653366),
    @(41592, 2.429152676597084),
    @(41774, 0.1648657400322122),
    @(41957, 1.153911995371715),
    @(42137, 0.9978838984052771),
    @(42321, 1.7958205029814),
    @(42503, -0.6421078521012049),
    @(42689, -0.4018974900578201),
    @(42867, 1.682996656587392),
    @(43053, 1.044515111522799),
    @(43145, 1.750939777838624),
    @(43235, -1.040686374073175),
    @(43326, -0.2631342424776193),
    @(43418, 0.8),
    @(43510, -0.8578072826671388),
    @(43600, 1),
    @(43691, 1.83469776744667),
    @(43783, -1.25257285039946),
    @(43875, 0.9704313550711277),
    @(43966, -3.1),
    @(44068, -20.33630385521235),
    @(44159, 18.14253003720006),
    @(44251, 4.544557841334466),
    @(44341, 1.849039957006937),
    @(44432, 0.5397909189180297),
    @(44525, -0.9815688985294599),
    @(44617, 4.784172743600905),
    @(44706, -2.12351799362105),
    @(44798, 0.2867860097325945),
    @(44890, 1.982042055103904),
    @(44981, -0.9572753870817365),
    @(45071, 0.4216256340467766),
    @(45163, -1.083093549703378),
    @(45254, -0.756372051364167),
    @(45345, -1.594471108675521),
    @(45436, 1.08522880488573),
    @(45534, -0.150193889177217),
    @(45618, -1.874498177919662),
    @(45713, -2.24094337053269),
    @(45800, 3.234504069499366),
    @(45891, -0.1262090121580286)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

$ws.Range("A54:A73").EntireRow.Delete()

Write-Host "done"